$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-10) are reordered so that the "Wall says La Loche school..."
# and "Court case for accused La Loche shooter adjourned" articles move to the
# top (rows 2 and 3), pushing the remaining rows down while keeping their
# original relative order. The last two rows stay where they were.

$ws.Cells.Item(2,1).Value2 = "Wall says La Loche school that was site of deadly shooting will be utilized"
$ws.Cells.Item(2,2).Value2 = "2016-02-09T19:31:00UTC"
$ws.Cells.Item(2,3).Value2 = 18
$ws.Cells.Item(2,4).Value2 = "day_2_to_30"
$ws.Cells.Item(2,5).Value2 = "http://saskatoon.ctvnews.ca/wall-says-la-loche-school-that-was-site-of-deadly-shooting-will-be-utilized-1.2771633"

$ws.Cells.Item(3,1).Value2 = "Court case for accused La Loche shooter adjourned"
$ws.Cells.Item(3,2).Value2 = "2016-07-19T06:00:00UTC"
$ws.Cells.Item(3,3).Value2 = 179
$ws.Cells.Item(3,4).Value2 = "day_31_beyond"
$ws.Cells.Item(3,5).Value2 = "http://battlefordsnow.com/article/525319/court-case-accused-la-loche-shooter-adjourned"

$ws.Cells.Item(4,1).Value2 = "Canada shootings: Four killed in Saskatchewan"
$ws.Cells.Item(4,2).Value2 = "2016-01-23T08:36:34UTC"
$ws.Cells.Item(4,3).Value2 = 1
$ws.Cells.Item(4,4).Value2 = "day_1"
$ws.Cells.Item(4,5).Value2 = "https://www.bbc.com/news/world-us-canada-35388594?ns_mchannel=social&ns_campaign=bbc_breaking&ns_source=twitter&ns_linkname=news_central"

$ws.Cells.Item(5,1).Value2 = "Trudeau arrives in La Loche, tells community to stand together after shooting"
$ws.Cells.Item(5,2).Value2 = "2016-01-29T11:17:09UTC"
$ws.Cells.Item(5,3).Value2 = 7
$ws.Cells.Item(5,4).Value2 = "day_2_to_30"
$ws.Cells.Item(5,5).Value2 = "http://globalnews.ca/news/2485895/trudeau-wall-head-to-la-loche-teachers-write-letter-to-students/"

$ws.Cells.Item(6,1).Value2 = "Teen charged in deadly La Loche, Sask., shooting pleads guilty"
$ws.Cells.Item(6,2).Value2 = "2016-10-28T16:34:00UTC"
$ws.Cells.Item(6,3).Value2 = 280
$ws.Cells.Item(6,4).Value2 = "day_31_beyond"
$ws.Cells.Item(6,5).Value2 = "http://www.cbc.ca/news/canada/saskatoon/court-appearance-la-loche-shooting-oct28-1.3825928"

$ws.Cells.Item(7,1).Value2 = "Saskatchewan will reach out to U.S. officials after La Loche shooting: Wall"
$ws.Cells.Item(7,2).Value2 = "2016-01-24T08:29:00UTC"
$ws.Cells.Item(7,3).Value2 = 2
$ws.Cells.Item(7,4).Value2 = "day_2_to_30"
$ws.Cells.Item(7,5).Value2 = "http://www.ctvnews.ca/canada/saskatchewan-will-reach-out-to-u-s-officials-after-la-loche-shooting-wall-1.2749780"

$ws.Cells.Item(8,1).Value2 = "La Loche Community School to remain closed"
$ws.Cells.Item(8,2).Value2 = "2016-02-03T00:00:00UTC"
$ws.Cells.Item(8,3).Value2 = 12
$ws.Cells.Item(8,4).Value2 = "day_2_to_30"
$ws.Cells.Item(8,5).Value2 = "http://www.cbc.ca/news/canada/saskatoon/la-loche-community-school-to-remain-closed-1.3432711"

$ws.Cells.Item(9,1).Value2 = "Canada gunman kills four in worst shooting since ’89"
$ws.Cells.Item(9,2).Value2 = "2016-01-23T18:23:57UTC"
$ws.Cells.Item(9,3).Value2 = 1
$ws.Cells.Item(9,4).Value2 = "day_1"
$ws.Cells.Item(9,5).Value2 = "http://www.japantimes.co.jp/news/2016/01/23/world/crime-legal-world/five-dead-two-critical-shootings-northern-saskatchewan/#.VqSyHxSJ_Z5"

$ws.Cells.Item(10,1).Value2 = "3 La Loche shooting victims released from hospital"
$ws.Cells.Item(10,2).Value2 = "1970-01-01T00:00:00UTC"
$ws.Cells.Item(10,3).Value2 = "unknown"
$ws.Cells.Item(10,4).Value2 = "unknown"
$ws.Cells.Item(10,5).Value2 = "http://www.cbc.ca/news/canada/saskatoon/la-loche-shootings-three-victims-released-from-hospital-1.3422108"
